# Rename the diff-table column headers in row 1 so the "_old"/"_new" suffixes
# become "_FV2410"/"_FV2504" (matching the respective input file names).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

# Columns A-J (1-10) carry the "_old" -> "_FV2410" headers.
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $baseNames[$i] + "_FV2410"
}

# Column K (11) is "diff" - unchanged.

# Columns L-U (12-21) carry the "_new" -> "_FV2504" headers.
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $baseNames[$i] + "_FV2504"
}

# Turn the data range into a proper Excel Table ("Table1").
$tableRange = $ws.Range("A1:U72")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# Freeze the header row.
$ws.Range("A2").Select() | Out-Null
[void]($excel.ActiveWindow.FreezePanes = $true)

Write-Output "edit complete"
